# Weekly update: a new Ciboulette price record (week of 2022-08-19,
# serial date 44792) is inserted at the top of the Femacal de La Calera
# series (row 153), pushing every existing record down by one row and
# appending the formerly-last record (2021-07-05, serial 44382) as the
# new final row (355).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 153:354 down to 154:355, carrying every column (values +
# formatting) with them - this is what creates the new row 355 and
# leaves row 153 open for the new observation.
$ws.Rows.Item(153).Insert()

# Populate the newly-opened row 153 with the new observation.
$ws.Cells.Item(153, 1).Value = 3
$ws.Cells.Item(153, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(153, 3).Value = "Coquimbo"
$ws.Cells.Item(153, 4).Value = 44792
$ws.Cells.Item(153, 5).Value = 5
$ws.Cells.Item(153, 6).Value = 100112039
$ws.Cells.Item(153, 7).Value = "Ciboulette"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 120
$ws.Cells.Item(153, 11).Value = 1500
$ws.Cells.Item(153, 12).Value = 1500
$ws.Cells.Item(153, 13).Value = 1500
$ws.Cells.Item(153, 14).Value = "`$/docena de atados"
$ws.Cells.Item(153, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(153, 16).Value = 500
$ws.Cells.Item(153, 17).Value = 3
$ws.Cells.Item(153, 18).Value = "Hortaliza"

# Match the date formatting used by the rest of column D.
$ws.Cells.Item(153, 4).NumberFormat = $ws.Cells.Item(154, 4).NumberFormat
